# chore: adapt column header formatting to respective input file names
#
# Rename header suffixes "_old" -> "_FV2310" and "_new" -> "_FV2404" across
# the header row (A1:U1), then wrap the sheet's used range (A1:U71) in an
# Excel Table ("Table1"), and finally freeze the header row (top row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Preserve the header row's existing cell formatting (bold / fill / border /
# centered / wrap) across the table creation below: stash a copy in a scratch
# row, strip direct formatting from the header so ListObjects.Add doesn't
# synthesize a header-row dxf override, then restore the stashed formatting.
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A200:U200")
$headerRange.Copy($scratch)
$headerRange.ClearFormats()

$tableRange = $ws.Range("A1:U71")
$listObject = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$listObject.Name = "Table1"
$listObject.TableStyle = $null

$scratch.Copy()
$headerRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$scratch.Clear()

# Freeze the top (header) row.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
